$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new status row (A4/B4) - "Hold/UnHold button - Just added"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Hold/UnHold button - Just added"

# The longer text in column B requires a wider column
$ws.Columns("B").ColumnWidth = 30.7109375

# The saved file's last active selection moved to G11
$ws.Activate()
$ws.Range("G11").Select()

$wb.Save()
